$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 (new last row) first inherits the old row 5 "last row" style (thicker bottom border) ---
$ws.Range("B5:G5").Copy()
$ws.Range("B26:G26").PasteSpecial(-4122)

# --- Rows 5 through 25 now take on the plain "body row" style used by rows 2-4 ---
$ws.Range("B2:G2").Copy()
$ws.Range("B5:G25").PasteSpecial(-4122)

# --- Fill rows 5 through 26 with their final data ---
$ws.Cells.Item(5,2).Value = 10333354
$ws.Cells.Item(5,3).Value = 'ولي العهد'
$ws.Cells.Item(5,4).Value = 'حدائق القبة'
$ws.Cells.Item(5,5).Value = 'منة الله مجدي'
$ws.Cells.Item(5,6).Value = 1026846805
$ws.Cells.Item(5,7).Value = 'Mennat Allah Magdy'
$ws.Cells.Item(6,2).Value = 10333396
$ws.Cells.Item(6,3).Value = 'المحكمة'
$ws.Cells.Item(6,4).Value = 'الزيتون و مصر الجديدة'
$ws.Cells.Item(6,5).Value = 'مصطفى علاء'
$ws.Cells.Item(6,6).Value = 1060459705
$ws.Cells.Item(6,7).Value = 'Mostafa Alaa'
$ws.Cells.Item(7,2).Value = 10331412
$ws.Cells.Item(7,3).Value = 'الخلفاوي'
$ws.Cells.Item(7,4).Value = 'شبرا'
$ws.Cells.Item(7,5).Value = 'جون مدحت'
$ws.Cells.Item(7,6).Value = 1281804393
$ws.Cells.Item(7,7).Value = 'John Medhat'
$ws.Cells.Item(8,2).Value = 10331413
$ws.Cells.Item(8,3).Value = 'كوبري الدقي'
$ws.Cells.Item(8,4).Value = 'المهندسين'
$ws.Cells.Item(8,5).Value = 'احمد نبيل'
$ws.Cells.Item(8,6).Value = '1225945969 / 01119986050'
$ws.Cells.Item(8,7).Value = 'Ahmed Nabil'
$ws.Cells.Item(9,2).Value = 10331414
$ws.Cells.Item(9,3).Value = 'بوابة 1'
$ws.Cells.Item(9,4).Value = 'و - مدينتي'
$ws.Cells.Item(9,5).Value = 'طارق عاطف'
$ws.Cells.Item(9,6).Value = 1117144455
$ws.Cells.Item(9,7).Value = 'Tarek Atef'
$ws.Cells.Item(10,2).Value = 10322712
$ws.Cells.Item(10,3).Value = 'المريوطيه'
$ws.Cells.Item(10,4).Value = 'فيصل'
$ws.Cells.Item(10,5).Value = 'مينا كارل'
$ws.Cells.Item(10,6).Value = 1152621423
$ws.Cells.Item(10,7).Value = 'Mina Carl'
$ws.Cells.Item(11,2).Value = 10331426
$ws.Cells.Item(11,3).Value = 'كارفور العبور'
$ws.Cells.Item(11,4).Value = 'العبور'
$ws.Cells.Item(11,5).Value = 'خالد وليد'
$ws.Cells.Item(11,6).Value = 1124220159
$ws.Cells.Item(11,7).Value = 'Khaled Waleed'
$ws.Cells.Item(12,2).Value = 10331422
$ws.Cells.Item(12,3).Value = 'الجراج'
$ws.Cells.Item(12,4).Value = 'الزيتون و مصر الجديدة'
$ws.Cells.Item(12,5).Value = 'روان سامح'
$ws.Cells.Item(12,6).Value = 1110071185
$ws.Cells.Item(12,7).Value = 'Rawan Sameh'
$ws.Cells.Item(13,2).Value = 10331421
$ws.Cells.Item(13,3).Value = 'اكاديمية السادات'
$ws.Cells.Item(13,4).Value = 'ي - المعادي'
$ws.Cells.Item(13,5).Value = 'شريف سيد'
$ws.Cells.Item(13,6).Value = 1013030118
$ws.Cells.Item(13,7).Value = 'Sherif Sayed'
$ws.Cells.Item(14,2).Value = 10333437
$ws.Cells.Item(14,3).Value = 'النافورة'
$ws.Cells.Item(14,4).Value = 'م - المقطم'
$ws.Cells.Item(14,5).Value = 'عفراء حاتم'
$ws.Cells.Item(14,6).Value = '1121009457 // 249964844267 WPP'
$ws.Cells.Item(14,7).Value = 'Afraa Hatim'
$ws.Cells.Item(15,2).Value = 10333354
$ws.Cells.Item(15,3).Value = 'ولي العهد'
$ws.Cells.Item(15,4).Value = 'حدائق القبة'
$ws.Cells.Item(15,5).Value = 'منة الله مجدي'
$ws.Cells.Item(15,6).Value = 1026846805
$ws.Cells.Item(15,7).Value = 'Mennat Allah Magdy'
$ws.Cells.Item(16,2).Value = 10333396
$ws.Cells.Item(16,3).Value = 'ميدان المحكمة'
$ws.Cells.Item(16,4).Value = 'الزيتون و مصر الجديدة'
$ws.Cells.Item(16,5).Value = 'مصطفى علاء'
$ws.Cells.Item(16,6).Value = 1060459705
$ws.Cells.Item(16,7).Value = 'Mostafa Alaa'
$ws.Cells.Item(17,2).Value = 10333385
$ws.Cells.Item(17,3).Value = 'بوابة 1'
$ws.Cells.Item(17,4).Value = 'و - مدينتي'
$ws.Cells.Item(17,5).Value = 'عبدالرحمن سليمان'
$ws.Cells.Item(17,6).Value = 1025613339
$ws.Cells.Item(17,7).Value = 'Abdulrahman Suliman'
$ws.Cells.Item(18,2).Value = 10333416
$ws.Cells.Item(18,3).Value = 'قسم الوايلي'
$ws.Cells.Item(18,4).Value = 'العباسية و الضاهر'
$ws.Cells.Item(18,5).Value = 'بسنت أحمد سمير'
$ws.Cells.Item(18,6).Value = 1270058439
$ws.Cells.Item(18,7).Value = 'Passant Ahmed Samir'
$ws.Cells.Item(19,2).Value = 10333417
$ws.Cells.Item(19,3).Value = 'كارفور العبور'
$ws.Cells.Item(19,4).Value = 'العبور'
$ws.Cells.Item(19,5).Value = 'مريم ابو بكر'
$ws.Cells.Item(19,6).Value = 1019687588
$ws.Cells.Item(19,7).Value = 'Mariam Abu Bakr'
$ws.Cells.Item(20,2).Value = 10333426
$ws.Cells.Item(20,3).Value = 'سلم البارون'
$ws.Cells.Item(20,4).Value = 'حلوان و زهراء المعادي'
$ws.Cells.Item(20,5).Value = 'مهاب علاء'
$ws.Cells.Item(20,6).Value = 1097696568
$ws.Cells.Item(20,7).Value = 'Mohab Alaa'
$ws.Cells.Item(21,2).Value = 10333432
$ws.Cells.Item(21,3).Value = 'مستشفى احمد ماهر'
$ws.Cells.Item(21,4).Value = 'التحرير'
$ws.Cells.Item(21,5).Value = 'مصطفى محمد'
$ws.Cells.Item(21,6).Value = 1091247732
$ws.Cells.Item(21,7).Value = 'Mustafa Mohamed'
$ws.Cells.Item(22,2).Value = 10333433
$ws.Cells.Item(22,3).Value = 'بهتيم'
$ws.Cells.Item(22,4).Value = 'دائري'
$ws.Cells.Item(22,5).Value = 'آية علي'
$ws.Cells.Item(22,6).Value = 1004788915
$ws.Cells.Item(22,7).Value = 'Aya Ali'
$ws.Cells.Item(23,2).Value = 10333434
$ws.Cells.Item(23,3).Value = 'كارفور العبور'
$ws.Cells.Item(23,4).Value = 'العبور'
$ws.Cells.Item(23,5).Value = 'احمد مجدي'
$ws.Cells.Item(23,6).Value = 1032409151
$ws.Cells.Item(23,7).Value = 'Ahmed Magdy'
$ws.Cells.Item(24,2).Value = 10333435
$ws.Cells.Item(24,3).Value = 'مسجد السلام'
$ws.Cells.Item(24,4).Value = 'مدينة نصر'
$ws.Cells.Item(24,5).Value = 'ماهر علي دهب'
$ws.Cells.Item(24,6).Value = '1212874368 / 967777756816'
$ws.Cells.Item(24,7).Value = 'Maher Ali Dahab'
$ws.Cells.Item(25,2).Value = 10333439
$ws.Cells.Item(25,3).Value = 'سلم البارون'
$ws.Cells.Item(25,4).Value = 'حلوان و زهراء المعادي'
$ws.Cells.Item(25,5).Value = 'ابراهيم بلة'
$ws.Cells.Item(25,6).Value = '1505852821 / 255749772862'
$ws.Cells.Item(25,7).Value = 'Ibrahim Balla'
$ws.Cells.Item(26,2).Value = 10333436
$ws.Cells.Item(26,3).Value = 'كشري هند الحي العاشر'
$ws.Cells.Item(26,4).Value = 'مدينة نصر'
$ws.Cells.Item(26,5).Value = 'عائشه محمد'
$ws.Cells.Item(26,6).Value = '1555175582 / 1107578879'
$ws.Cells.Item(26,7).Value = 'Aisha Mohammed'

# --- Row 26 gets the taller row height that used to belong to row 5 ---
$ws.Rows.Item(26).RowHeight = 15.25

# --- Row 5 reverts to the default row height since it is no longer the last row ---
$ws.Rows.Item(5).AutoFit()

# --- Column F gets wider to fit the new, longer phone-number strings ---
$ws.Columns.Item(6).ColumnWidth = 32

# --- Update the active selection ---
$ws.Range("D7").Select()
